$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.504.88'
$ws.Range('E2').Value = '  -0.50%  '
$ws.Range('D3').Value = '1.628.05'
$ws.Range('E3').Value = '  -0.05%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.82'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.13%  '
$ws.Range('E6').Value = '  +1.76%  '
$ws.Range('E7').Value = '  +0.25%  '
$ws.Range('E8').Value = '  -0.46%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0621'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.62%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.02'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.36%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0838'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.53%  '
$ws.Range('D12').Value = '1.855.33'
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').Value = '1.625.23'
$ws.Range('E13').Value = '  -0.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.11'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.72%  '
$ws.Range('E15').Value = '  -0.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.88'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.49%  '
$ws.Range('D17').Value = '26.533.39'
$ws.Range('E17').Value = '  -0.35%  '
$ws.Range('D18').Value = '0.0₃0739'
$ws.Range('E18').Value = '  +1.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '215.62'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.00'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E21').Value = '  -0.17%  '
$ws.Range('E22').Value = '  +1.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.31'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.98'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.05%  '
$ws.Range('E25').Value = '  +1.88%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.01'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.26%  '
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('E28').Value = '  +3.26%  '
$ws.Range('E29').Value = '  +1.30%  '
$ws.Range('E30').Value = '  -2.31%  '
$ws.Range('E31').Value = '  -1.25%  '
$ws.Range('E32').Value = '  +2.86%  '
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').Value = '1.220.97'
$ws.Range('E35').Value = '  +4.91%  '
$ws.Range('E36').Value = '  -1.63%  '
$ws.Range('E37').Value = '  +5.55%  '
$ws.Range('E38').Value = '  +0.24%  '
$ws.Range('E39').Value = '  -1.31%  '
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('E41').Value = '  -2.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.794'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.35'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.74%  '
$ws.Range('D44').Value = '1.765.05'
$ws.Range('E44').Value = '  +0.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.40'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.34%  '
$ws.Range('E46').Value = '  +1.92%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '54.83'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.70%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0511'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.11%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.65'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.55%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.409'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.18%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.01'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.43%  '
